# Highlight (yellow) the intro bullet + its first sub-bullet for
# "Exercício 1" on slide 11 and "Exercício 6" on slide 12.

$p = $ppt.ActivePresentation

# --- Slide 11: "Exercício 1: Verificar se um Número é Par" ---
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(2)
$tr11 = $shp11.TextFrame2.TextRange

# Paragraph 1: "Exercício 1: Verificar se um Número é Par "
$tr11.Paragraphs(1, 1).Font.Highlight.RGB = 65535

# Paragraph 2: "Crie uma função anônima que recebe um número como
# parâmetro e retorna `true` se o número for par, ou `false` se for ímpar. "
$tr11.Paragraphs(2, 1).Font.Highlight.RGB = 65535

# --- Slide 12: "Exercício 6: Calcular a Média de Três Números" ---
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(2)
$tr12 = $shp12.TextFrame2.TextRange

# Paragraph 1: "Exercício 6: Calcular a Média de Três Números"
$tr12.Paragraphs(1, 1).Font.Highlight.RGB = 65535

# Paragraph 2: "Crie uma função anônima que recebe três números como
# parâmetros e retorna a média desses números."
$tr12.Paragraphs(2, 1).Font.Highlight.RGB = 65535
